$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "Incident state changed to Closed"
#   -> "Incident state changed to " (run 1, unchanged formatting/rsid)
#      + "Completed" (run 2, new run, no rsid)
# ---------------------------------------------------------------------------

$rClosed = $d.Content
$rClosed.Find.Execute("Closed", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# Force a run split at this boundary by toggling a character formatting
# property on/off (mirrors how Word splits a run when applying/removing
# direct formatting to part of it), then assign the new text so the
# freshly created run carries no leftover rsid.
$rClosed.Bold = 1
$rClosed.Text = "Completed"
$rClosed.Bold = 0

# ---------------------------------------------------------------------------
# Edit 2: "System changes the status of incident to 'closed'."
#   -> "System changes the status of incident to '" (run 1, unchanged)
#      + "Completed" (run 2, new run, no rsid)
#      + "'" (run 3, new run, no rsid)
#      + "." (run 4, new run, no rsid)
# ---------------------------------------------------------------------------

$rWord = $d.Content
$rWord.Find.Execute("closed", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$wordStart = $rWord.Start

$rWord.Bold = 1
$rWord.Text = "Completed"
$rWord.Bold = 0

$newWordEnd = $wordStart + 9   # len("Completed")

$rQuote = $d.Range($newWordEnd, $newWordEnd + 1)
$rQuote.Bold = 1
$rQuote.Text = "ZZ"
$rQuote.Text = [string][char]0x2019
$rQuote.Bold = 0

$quoteEnd = $newWordEnd + 1

$rPeriod = $d.Range($quoteEnd, $quoteEnd + 1)
$rPeriod.Bold = 1
$rPeriod.Text = "ZZ"
$rPeriod.Text = "."
$rPeriod.Bold = 0
